$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (columns A:C), replacing the old rows 2-8 sample data
$data = @(
    @(32, 5, 90000),
    @(28, 3, 65000),
    @(45, 15, 150000),
    @(36, 7, 60000),
    @(52, 20, 200000),
    @(29, 2, 55000),
    @(42, 12, 120000),
    @(31, 4, 80000),
    @(26, 1, 45000),
    @(38, 10, 110000),
    @(29, 3, 75000),
    @(48, 18, 140000),
    @(35, 6, 65000),
    @(40, 14, 130000),
    @(27, 2, 40000)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Apply vertical-center alignment to the whole new data block (A2:C16) -
# this creates the new cellXfs entry with alignment vertical="center"
$ws.Range("A2:C16").VerticalAlignment = -4108

# Update selection to match the saved view state
$ws.Range("E6").Select() | Out-Null
